$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 2017
$ws.Range("B8").Value = "GL_alvar"
$ws.Range("C8").Value = 576.436
$ws.Range("D8").Value = 125.247

$ws.Range("A9").Value = 2017
$ws.Range("B9").Value = "MB_alvar"
$ws.Range("C9").Value = 93.6
$ws.Range("D9").Value = 64.409

$ws.Range("A10").Value = 2017
$ws.Range("B10").Value = "Prairie"
$ws.Range("C10").Value = 65.652
$ws.Range("D10").Value = 63.235

$ws.Range("D10").Select()
